$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the scraped cryptocurrency price/volume refresh.
# Rows 39 and 40 swap coin identity (Maker <-> VeChain) together with their
# price/volume figures; all other rows just get refreshed price/volume text.
# A leading apostrophe is prefixed to values that would otherwise be
# auto-parsed as numbers by Excel, so the cell keeps the exact original text
# (e.g. trailing zeros such as "0.6150") instead of silently becoming numeric.

$ws.Range("D2").Value = "29.071.71"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "1.835.63"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'242.67"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "'0.6150"
$ws.Range("E6").Value = "  -3.20%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.07476"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'0.2925"
$ws.Range("E9").Value = "  -0.92%  "
$ws.Range("D10").Value = "'23.10"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "'0.07688"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.834.39"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'4.996"
$ws.Range("E13").Value = "  -0.34%  "
$ws.Range("D14").Value = "'0.6718"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'82.66"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "'0.000009145"
$ws.Range("E16").Value = "  -4.34%  "
$ws.Range("D17").Value = "'5.909"
$ws.Range("E17").Value = "  -3.14%  "
$ws.Range("D18").Value = "29.042.22"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("D19").Value = "2.082.00"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'232.45"
$ws.Range("E20").Value = "  +1.94%  "
$ws.Range("E21").Value = "  -0.02%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'7.196"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'158.87"
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "'0.1397"
$ws.Range("E26").Value = "  -2.14%  "
$ws.Range("D27").Value = "'8.489"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("D29").Value = "'1.499"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").Value = "'4.152"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'4.115"
$ws.Range("D32").Value = "'0.05497"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("D33").Value = "'1.203"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'0.7380"
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("D37").Value = "'2.661"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("D38").Value = "'2.776"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01779"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.213.73"
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").Value = "'6.441"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("D42").Value = "'0.8934"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").Value = "'101.89"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").Value = "1.980.65"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'65.42"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "'0.5088"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "'0.00000000119"
$ws.Range("E48").Value = "  -5.79%  "
$ws.Range("D49").Value = "'0.4069"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'9.153"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +0.60%  "
